$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04562342498972654
$ws.Range("C2").Value = 1.152029883462799
$ws.Range("B3").Value = 0.06693906232308582
$ws.Range("C3").Value = 2.992764033751882
